$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "televisor"
$ws.Range("B2").Value = "10"

$ws.Range("D6").Select()
